# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45177 (2023-09-08) to 45178 (2023-09-09).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 329 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45178
